$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 8c3deefb-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-30 14:59:24"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 8c3deefb-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-30 14:59:19"
$wsZhCn.Range("K3").Value = "2016-08-30 14:59:53"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 8c3deefb-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-30 14:59:24"
$wsDeDe.Range("K3").Value = "2016-08-30 15:00:04"
